# rnaSample_468.xlsx -- fixed harvester column in rnaSamples
# Holly added "S.GISH" as a harvester value in bioSamples; propagate that
# fix to the "harvester" column (column B) of this rnaSamples sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the harvester column (B2:B14): Retrofitted_468 -> S.GISH -------
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
}

# --- Cosmetic clean-up that accompanied the edit in the source file -----

# The G1 header cell ("roboticRNAPrep") carried its own (visually
# identical) bold-font style, distinct from the rest of the header row.
# Nudge the font so the engine re-resolves it against the already-used
# header font/style instead of keeping its own redundant copy.
$ws.Range("G1").Font.Bold = $false
$ws.Range("G1").Font.Bold = $true

# Column B is now slightly wider than the rest of the sheet.
$ws.Columns.Item(2).ColumnWidth = 8.83 - 5/6

# The selection left active in the saved file moved to C4.
$ws.Range("C4").Select()
